$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [D value, E value] updates (price, volume/1h change)
$updates = @{
    2  = @("24.876.19", "  +0.67%  ")
    3  = @("1.661.99",  "  +0.74%  ")
    4  = @($null,       "  -0.61%  ")
    5  = @("325.56",    "  +5.65%  ")
    6  = @("1.000",     "  -0.27%  ")
    7  = @("0.3626",    "  -0.31%  ")
    8  = @("47.70",     "  +2.07%  ")
    9  = @("0.3265",    "  -0.86%  ")
    10 = @("1.132",     "  -0.15%  ")
    11 = @("0.07075",   "  -0.44%  ")
    12 = @("1.000",     "  -0.30%  ")
    13 = @("6.038",     "  -0.22%  ")
    14 = @("19.49",     "  -1.34%  ")
    15 = @("1.657.05",  "  +0.38%  ")
    16 = @("6.600",     "  -1.05%  ")
    17 = @("0.00001045","  -1.77%  ")
    18 = @("0.06599",   "  +0.05%  ")
    19 = @("0.9996",    "  -0.35%  ")
    20 = @("78.98",     "  -0.88%  ")
    21 = @($null,       "  -1.78%  ")
    22 = @("15.78",     "  -3.60%  ")
    23 = @("12.56",     "  +2.96%  ")
    24 = @("24.869.17", "  +0.68%  ")
    25 = @("2.450",     "  +1.59%  ")
    26 = @("2.426",     "  -4.56%  ")
    27 = @("148.58",    "  +0.04%  ")
    28 = @("18.64",     "  -3.19%  ")
    29 = @("1.836.54",  "  +0.05%  ")
    30 = @("125.19",    "  -2.16%  ")
    31 = @("1.187",     "  +4.94%  ")
    32 = @("4.083",     "  -1.23%  ")
    33 = @("5.721",     "  -6.69%  ")
    34 = @("0.08435",   "  -0.45%  ")
    35 = @("1.649",     "  -4.29%  ")
    36 = @("12.18",     "  -4.34%  ")
    37 = @("1.282",     "  +4.53%  ")
    38 = @("5.168",     "  -1.06%  ")
    39 = @("0.02269",   "  -0.92%  ")
    40 = @("0.06111",   "  -1.87%  ")
    43 = @("0.9998",    "  -0.21%  ")
    44 = @("0.5930",    "  -2.42%  ")
    47 = @("0.5622",    "  -2.68%  ")
    48 = @("124.85",    "  +1.07%  ")
    49 = @("1.948",     "  -2.24%  ")
    50 = @("0.06984",   "  -0.99%  ")
    51 = @("1.190",     "  +1.79%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals[0] -ne $null) {
        $ws.Range("D$row").Value = $vals[0]
    }
    $ws.Range("E$row").Value = $vals[1]
}

# Rows 41/42 swap content (Algorand <-> FraxShare) with new D/E values
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.2073"
$ws.Range("E41").Value = "  -1.35%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "8.292"
$ws.Range("E42").Value = "  -0.57%  "

# Rows 45/46 swap content (EnergySwap <-> PancakeSwap) with new D/E values
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "13.37"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.835"
$ws.Range("E46").Value = "  +2.11%  "
